$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F40").Value = 103
$ws.Range("G40").Value = 3216.69
$ws.Range("B71").Value = 54844.3
$ws.Range("F105").Value = 221
$ws.Range("G105").Value = 17938.57
$ws.Range("F126").Value = 55
$ws.Range("G126").Value = 13533.85
$ws.Range("F141").Value = 385
$ws.Range("G141").Value = 7492.1
$ws.Range("B143").Value = 278439.5
$ws.Range("F163").Value = 2
$ws.Range("G163").Value = 79.72
$ws.Range("B176").Value = 14867.26
$ws.Range("F195").Value = 61
$ws.Range("G195").Value = 2581.52
$ws.Range("B205").Value = 26677.53
$ws.Range("F235").Value = 27
$ws.Range("G235").Value = 2279.34
$ws.Range("F238").Value = 18
$ws.Range("G238").Value = 913.3200000000001
$ws.Range("F241").Value = 47
$ws.Range("G241").Value = 1856.03
$ws.Range("B250").Value = 27469.91
$ws.Range("F257").Value = 24
$ws.Range("G257").Value = 1800.48
$ws.Range("B262").Value = 15527.42
$ws.Range("F296").Value = 71
$ws.Range("G296").Value = 6022.93
$ws.Range("F298").Value = 95
$ws.Range("G298").Value = 6546.45
$ws.Range("B301").Value = 13241.6
$ws.Range("F333").Value = 40
$ws.Range("G333").Value = 1282.8
$ws.Range("F363").Value = 45
$ws.Range("G363").Value = 5300.55
$ws.Range("B382").Value = 126008.59
$ws.Range("F393").Value = 33
$ws.Range("G393").Value = 2737.02
$ws.Range("F405").Value = 107
$ws.Range("G405").Value = 14666.49
$ws.Range("F420").Value = 280
$ws.Range("G420").Value = 14596.4
$ws.Range("F425").Value = 1
$ws.Range("G425").Value = 85.72
$ws.Range("F428").Value = 4
$ws.Range("G428").Value = 497.68
$ws.Range("F430").Value = 14
$ws.Range("G430").Value = 1457.26
$ws.Range("F431").Value = 22
$ws.Range("G431").Value = 1953.6
$ws.Range("F434").Value = 171
$ws.Range("G434").Value = 17288.1
$ws.Range("F436").Value = 59
$ws.Range("G436").Value = 11716.81
$ws.Range("F442").Value = 163
$ws.Range("G442").Value = 9638.190000000001
$ws.Range("F452").Value = 61
$ws.Range("G452").Value = 12185.97
$ws.Range("F453").Value = 32
$ws.Range("G453").Value = 5951.68
$ws.Range("F454").Value = 258
$ws.Range("G454").Value = 18139.98
$ws.Range("F456").Value = 21
$ws.Range("G456").Value = 3174.57
$ws.Range("F461").Value = 59
$ws.Range("G461").Value = 18910.68
$ws.Range("B467").Value = 411858.42
$ws.Range("F492").Value = 3
$ws.Range("G492").Value = 108.3
$ws.Range("F496").Value = 18
$ws.Range("G496").Value = 367.02
$ws.Range("B500").Value = 1394.63
$ws.Range("F559").Value = 101
$ws.Range("G559").Value = 4787.4
$ws.Range("B571").Value = 53485.78
$ws.Range("F574").Value = 7
$ws.Range("G574").Value = 353.85
$ws.Range("F576").Value = 0
$ws.Range("G576").Value = 0
$ws.Range("B582").Value = 1294.75
$ws.Range("F587").Value = 648
$ws.Range("G587").Value = 8715.6
$ws.Range("F588").Value = 534
$ws.Range("G588").Value = 7022.1
$ws.Range("F589").Value = 634
$ws.Range("G589").Value = 8121.54
$ws.Range("F591").Value = 225
$ws.Range("G591").Value = 4439.25
$ws.Range("F592").Value = 352
$ws.Range("G592").Value = 5783.36
$ws.Range("F596").Value = 463
$ws.Range("G596").Value = 7509.86
$ws.Range("F597").Value = 148
$ws.Range("G597").Value = 2880.08
$ws.Range("F598").Value = 1048
$ws.Range("G598").Value = 6895.84
$ws.Range("F602").Value = 374
$ws.Range("G602").Value = 9836.200000000001
$ws.Range("F603").Value = 312
$ws.Range("G603").Value = 5126.16
$ws.Range("B605").Value = 110801.08
$ws.Range("F610").Value = 3
$ws.Range("G610").Value = 651.24
$ws.Range("B620").Value = 10662.34
$ws.Range("F669").Value = 684
$ws.Range("G669").Value = 4651.2
$ws.Range("F671").Value = 368
$ws.Range("G671").Value = 4743.52
$ws.Range("F672").Value = 662
$ws.Range("G672").Value = 13140.7
$ws.Range("F674").Value = 325
$ws.Range("G674").Value = 5372.25
$ws.Range("B677").Value = 43087.16
$ws.Range("F724").Value = 23
$ws.Range("G724").Value = 508.53
$ws.Range("B739").Value = 8135.42
$ws.Range("F743").Value = 41
$ws.Range("G743").Value = 4325.91
$ws.Range("F757").Value = 28
$ws.Range("G757").Value = 2059.4
$ws.Range("F758").Value = 79
$ws.Range("G758").Value = 9703.57
$ws.Range("B761").Value = 54300.02
$ws.Range("F767").Value = 109
$ws.Range("G767").Value = 2964.8
$ws.Range("B770").Value = 68895.12
$ws.Range("F796").Value = 281
$ws.Range("G796").Value = 4456.66
$ws.Range("F800").Value = 108
$ws.Range("G800").Value = 4663.44
$ws.Range("F802").Value = 183
$ws.Range("G802").Value = 7901.94
$ws.Range("F803").Value = 35
$ws.Range("G803").Value = 1531.6
$ws.Range("B804").Value = 36713.88
$ws.Range("F838").Value = 14
$ws.Range("G838").Value = 1125.32
$ws.Range("B839").Value = 4174.07
$ws.Range("F877").Value = 14
$ws.Range("G877").Value = 3436.58
$ws.Range("B880").Value = 27488.38
$ws.Range("F886").Value = 75
$ws.Range("G886").Value = 8365.5
$ws.Range("B901").Value = 51055.19
$ws.Range("F919").Value = 1
$ws.Range("G919").Value = 59.94
$ws.Range("F923").Value = 81
$ws.Range("G923").Value = 8335.709999999999
$ws.Range("F925").Value = 8
$ws.Range("G925").Value = 555.04
$ws.Range("F926").Value = 32
$ws.Range("G926").Value = 1178.56
$ws.Range("B933").Value = 36963.2
$ws.Range("F936").Value = 67
$ws.Range("G936").Value = 2505.8
$ws.Range("F939").Value = 149
$ws.Range("G939").Value = 5572.6
$ws.Range("F941").Value = 128
$ws.Range("G941").Value = 4787.2
$ws.Range("B942").Value = 13841.7
$ws.Range("F987").Value = 47
$ws.Range("G987").Value = 1818.43
$ws.Range("B997").Value = 3520.31
$ws.Range("F999").Value = 1411
$ws.Range("G999").Value = 230148.21
$ws.Range("F1002").Value = 69
$ws.Range("G1002").Value = 10202.34
$ws.Range("F1003").Value = 181
$ws.Range("G1003").Value = 12217.5
$ws.Range("B1005").Value = 271107.41
$ws.Range("B1012").Value = 2628904.3
$ws.Range("B1013").Value = 2628904.3
